# prepa_certif.xlsx revision:
#  - "Sprint backlog" row moves up (was row 18 -> row 14)
#  - new "Sprint planning" / "La scrum team peut y assister." row inserted at row 16
#  - remaining rows (Daily Scrum, Product backlog, DONE, scrum team, Scrum artifacts)
#    shift down to make room
#  - "scrum team" row gets an extra sentence appended to its description
#  - selection moves to B16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Stage the rows that need to move, in a scratch area far below the data,
#    so the move can happen even though source/destination ranges overlap.
# ---------------------------------------------------------------------------
$ws.Range("A14:B14").Copy($ws.Range("A200"))   # Daily Scrum / dev team
$ws.Range("A16:B16").Copy($ws.Range("A202"))   # Product backlog / (rich text)
$ws.Range("A18:B18").Copy($ws.Range("A204"))   # Sprint backlog / changes during sprint
$ws.Range("A20:B20").Copy($ws.Range("A206"))   # DONE / dev team (responsable...)
$ws.Range("A22:B22").Copy($ws.Range("A208"))   # scrum team / différentes équipes...
$ws.Range("A24:B24").Copy($ws.Range("A210"))   # Scrum artifacts / sprint backlog,...

# ---------------------------------------------------------------------------
# 2. Wipe the area that is being reshuffled (rows 14-24) clean.
# ---------------------------------------------------------------------------
$ws.Range("A14:B24").Clear()

# old spacer rows 17 & 19 no longer exist in the new layout
$ws.Range("A17:B19").Clear()
$ws.Rows("17:19").AutoFit()

# ---------------------------------------------------------------------------
# 3. Write the new layout back, row by row.
# ---------------------------------------------------------------------------

# Row 13 : blank spacer (label-style only, column A)
$ws.Range("A13").Font.Bold = $true
$ws.Range("A13").Font.Size = 14
$ws.Range("A13").HorizontalAlignment = -4131
$ws.Range("A13").VerticalAlignment = -4108
$ws.Range("A13").WrapText = $true

# Row 24 : scrum team (was row 22) - description gains an extra sentence.
# Written first so the new shared string ("...Elle fait le SPRINT PLANNING.")
# lands in the same shared-strings slot order as the source file.
$ws.Range("A208:B208").Copy($ws.Range("A24"))
$cur = $ws.Range("B24").Value2
$ws.Range("B24").Value2 = $cur + "`nElle fait le SPRINT PLANNING."
$ws.Rows("24:24").RowHeight = 43.5

# Row 16 : brand-new "Sprint planning" row
$ws.Range("A16").Value2 = "Sprint planning"
$ws.Range("B16").Value2 = "La scrum team peut y assister."
$ws.Range("A16").Font.Bold = $true
$ws.Range("A16").Font.Size = 14
$ws.Range("A16").HorizontalAlignment = -4131
$ws.Range("A16").VerticalAlignment = -4108
$ws.Range("A16").WrapText = $true
$ws.Range("B16").Font.Bold = $false
$ws.Range("B16").Font.Size = 11
$ws.Range("B16").HorizontalAlignment = -4131
$ws.Range("B16").VerticalAlignment = -4108
$ws.Range("B16").WrapText = $true
$ws.Rows("16:16").RowHeight = 18

# Row 14 : Sprint backlog (was row 18)
$ws.Range("A204:B204").Copy($ws.Range("A14"))
$ws.Rows("14:14").RowHeight = 14

# Row 15 : blank spacer (label-style only, column A)
$ws.Range("A15").Font.Bold = $true
$ws.Range("A15").Font.Size = 14
$ws.Range("A15").HorizontalAlignment = -4131
$ws.Range("A15").VerticalAlignment = -4108
$ws.Range("A15").WrapText = $true
$ws.Rows("15:15").RowHeight = 14

# Row 18 : Daily Scrum (was row 14)
$ws.Range("A200:B200").Copy($ws.Range("A18"))

# Row 20 : Product backlog (was row 16, keeps its rich-text formatting)
$ws.Range("A202:B202").Copy($ws.Range("A20"))
$ws.Rows("20:20").RowHeight = 58

# Row 21 : blank spacer (no cells)
$ws.Rows("21:21").RowHeight = 14

# Row 22 : DONE (was row 20)
$ws.Range("A206:B206").Copy($ws.Range("A22"))

# Row 23 : blank spacer (no cells)
$ws.Rows("23:23").RowHeight = 14

# Row 25 : blank spacer (no cells)
$ws.Rows("25:25").RowHeight = 14.5

# Row 26 : Scrum artifacts (was row 24)
$ws.Range("A210:B210").Copy($ws.Range("A26"))

# ---------------------------------------------------------------------------
# 4. Clean up the scratch area used for staging the moved rows.
# ---------------------------------------------------------------------------
$ws.Range("A200:B210").Clear()
$ws.Rows("200:210").AutoFit()

# ---------------------------------------------------------------------------
# 5. Restore the selection to match the saved file.
# ---------------------------------------------------------------------------
$ws.Range("B16").Select()
